# Update "想去人数" (want-to-go count) / "最低票价" (lowest price) figures
# on the "展览" and "全部类型" sheets to the freshly-scraped values
# (gh-pages output regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" -------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")

$ws.Range("F3").Value  = 5476
$ws.Range("F4").Value  = 14
$ws.Range("F5").Value  = 7594
$ws.Range("F6").Value  = 49
$ws.Range("F9").Value  = 607
$ws.Range("F11").Value = 35
$ws.Range("F12").Value = 4382
$ws.Range("F16").Value = 2968
$ws.Range("F17").Value = 580
$ws.Range("F20").Value = 539
$ws.Range("F21").Value = 466
$ws.Range("F22").Value = 474
$ws.Range("F23").Value = 333
$ws.Range("F24").Value = 114
$ws.Range("F25").Value = 1716
$ws.Range("F26").Value = 1231
$ws.Range("G26").Value = 45
$ws.Range("F27").Value = 99
$ws.Range("F28").Value = 1423
$ws.Range("F30").Value = 591
$ws.Range("F34").Value = 0
$ws.Range("F35").Value = 68
$ws.Range("F38").Value = 3037
$ws.Range("F39").Value = 713
$ws.Range("F41").Value = 129
$ws.Range("F43").Value = 657

# ---- Sheet "全部类型" ----------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")

$ws.Range("F3").Value  = 5476
$ws.Range("F4").Value  = 14
$ws.Range("F5").Value  = 7594
$ws.Range("F6").Value  = 49
$ws.Range("F9").Value  = 607
$ws.Range("F11").Value = 35
$ws.Range("F12").Value = 4382
$ws.Range("F16").Value = 2968
$ws.Range("F17").Value = 580
$ws.Range("F20").Value = 539
$ws.Range("F21").Value = 466
$ws.Range("F22").Value = 474
$ws.Range("F24").Value = 333
$ws.Range("F25").Value = 114
$ws.Range("F26").Value = 1716
$ws.Range("F27").Value = 1231
$ws.Range("G27").Value = 45
$ws.Range("F28").Value = 99
$ws.Range("F29").Value = 1423
$ws.Range("F31").Value = 591
$ws.Range("F36").Value = 68
$ws.Range("F39").Value = 3037
$ws.Range("F41").Value = 713
$ws.Range("F43").Value = 129
$ws.Range("F45").Value = 658
